$wb = $excel.ActiveWorkbook
$cms = $wb.Worksheets.Item("CMS")
$cms.Copy([System.Reflection.Missing]::Value, $cms)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "GS"

# clear existing content/formatting from the copied sheet
$ws.Cells.Clear()
$ws.Rows.Item(1).ClearFormats()
$ws.Rows.Item(1).EntireRow.AutoFit()

$headers = @("Contact_ID", "Contact_Date", "Contact_Type_Code", "Contact_Type_Desc", "OM_Name", "OM_Key", "OM_Grade", "OM_Team_Key", "OM_Provider_Code")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
}
$headerRange = $ws.Range("A1:I1")
$headerRange.Font.Color = 0
[void]$ws.Range("A1:I1").Select()
